$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCol = $ws.Range("A76:A77")

# Column A stores the date labels as plain text, not real dates - force
# text formatting so Excel doesn't coerce the strings into date serials
# while we (re)write them.
$dateCol.NumberFormat = "@"

# Fix the mislabeled date on the existing last row (2020/09/03 -> 2021/09/03)
$ws.Range("A76").Value = "2021/09/03"

# Append the new row from DGS's 2021/09/06 report
$ws.Range("A77").Value = "2021/09/06"
$ws.Range("B77").Value = 276
$ws.Range("C77").Value = 283.8
$ws.Range("D77").Value = 0.92
$ws.Range("E77").Value = 0.93

# Restore the original date-display number format so column A keeps the
# same look/style as the rest of the table.
$dateCol.NumberFormat = "yyyy/mm/dd"

$ws.Range("A77").Select()
